$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Replace the mixed-formatting "API design and programming" text in C7
# with the plain string "Restful API programming".
$ws.Range("C7").Value = "Restful API programming"
